$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 to include the additional assistant/user turns
$ws.Range("A2").Value = '[{"role": "assistant", "content": "Chào cậu! Bữa nay thầy giáo robot trên sao HỎA dạy tớ cách luyện nói theo nhịp như… hát rap đó! Mỗi cụm từ giống như một đoạn nhạc nhỏ, có nhịp điệu riêng."}, {"role": "assistant", "content": "Cứ theo nhịp là nói được tiếng Anh hay hơn liền! Giờ tụi mình thử cùng luyện nha. Cậu đã sẵn sàng chưa?"}, {"role": "user", "content": "[NEXT]"}]'

# Update C2 to the new content (matches the first assistant message text per the target diff)
$ws.Range("C2").Value = "Chào cậu! Bữa nay thầy giáo robot trên sao HỎA dạy tớ cách luyện nói theo nhịp như… hát rap đó! Mỗi cụm từ giống như một đoạn nhạc nhỏ, có nhịp điệu riêng."

# Remove rows 3 through 7 entirely (they are no longer part of the data)
$ws.Range("A3:C7").EntireRow.Delete()
